# Y3_B2526_GIT_&_Liver_schedule.xlsx
# Bumps the "Duration" (column G) values for a batch of sessions from 75 to
# 90 minutes, and one session (row 32) from 75 to 120 minutes (picking up
# the "no-fill" banding style used by odd data rows instead of its own
# even-row grey-fill style). Also updates the current selection to the
# header row (A1:XFD1), matching the saved sheet view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Selection saved in the sheet view (whole of row 1 selected).
[void]$ws.Range("A1:XFD1").Select()

# Row 32 is special-cased: duration becomes 120 AND it picks up the
# "no fill" formatting (the style already used a few rows down, e.g. G31)
# instead of keeping its own grey-fill style.
$ws.Range("G31").Copy()
$ws.Range("G32").PasteSpecial([int]-4122)  # xlPasteFormats
$ws.Range("G32").Value = 120

# Remaining rows: duration 75 -> 90, formatting/style unchanged.
$rows90 = 21, 29, 30, 31, 52, 60, 61, 62, 83, 91, 92, 93, 114, 122, 123, 124, `
          145, 153, 154, 155, 176, 184, 185, 186
foreach ($r in $rows90) {
    $ws.Range("G$r").Value = 90
}
